# Fruta / hortaliza, semanal
#
# A new weekly observation was inserted as the 2nd data row (worksheet row 3),
# pushing every row that used to be at position 3..91 down by one (to 4..92).
# Excel's native "insert row" operation reproduces exactly that shift
# (including carrying the date-column number format down to the new row),
# so we use it here instead of rewriting every row by hand.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 3 (old rows 3..91 -> 4..92).
$ws.Rows.Item(3).Insert()

# Populate the new row 3 with the inserted observation.
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C3").Value = "Arica y Parinacota"
$ws.Range("D3").Value = 44643
$ws.Range("E3").Value = 15
$ws.Range("F3").Value = "Fruta"
$ws.Range("G3").Value = 100102
$ws.Range("H3").Value = "Cítricos"
$ws.Range("I3").Value = 100102004
$ws.Range("J3").Value = "Mandarina"
$ws.Range("K3").Value = "Murcott"
$ws.Range("L3").Value = "Tercera"
$ws.Range("M3").Value = 250
$ws.Range("N3").Value = 13000
$ws.Range("O3").Value = 14000
$ws.Range("P3").Value = 13500
$ws.Range("Q3").Value = "$/caja 20 kilos"
$ws.Range("R3").Value = "Región de Coquimbo"
$ws.Range("S3").Value = 675
$ws.Range("T3").Value = 20
